# Apply updated dSF (column F) values as part of a data repull / recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = -3
    10 = 1
    11 = -2
    12 = -4
    30 = -4
    31 = -2
    32 = 3
    34 = 6
    36 = 1
    43 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
